# Swap the order of names in the "Recorded By" column (G) for every row
# where the value is exactly "System, dnasr281@gmail.com" or
# "dnasr281@gmail.com, System" - flipping it to the other ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$valA = "System, dnasr281@gmail.com"
$valB = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $valA) {
        $cell.Value2 = $valB
    } elseif ($val -eq $valB) {
        $cell.Value2 = $valA
    }
}
